$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "pontos notáveis - incremento na tabela de ranking": the share columns
# (E and F, rows 2-7) were stored as raw fractions (0-1) but need to be on
# the 0-100 percentage scale, so bump each value by a factor of 100 while
# leaving every other cell (and the existing 0.00% number format) untouched.
$range = $ws.Range("E2:F7")
foreach ($cell in $range.Cells) {
    $cell.Value = $cell.Value2 * 100
}
